$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.02"
$ws.Range("E2").Value = "'-0.57%"
$ws.Range("D3").Value = "'29.18"
$ws.Range("E3").Value = "'-1.70%"
$ws.Range("D4").Value = "'5.260"
$ws.Range("E4").Value = "'1.79%"
$ws.Range("D5").Value = "'0.05707"
$ws.Range("E5").Value = "'0.01%"
$ws.Range("E6").Value = "'0.23%"
$ws.Range("D7").Value = "'3.192"
$ws.Range("E7").Value = "'3.91%"
$ws.Range("D8").Value = "'0.8510"
$ws.Range("E8").Value = "'-0.84%"
$ws.Range("D9").Value = "'0.8586"
$ws.Range("E9").Value = "'-1.10%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1371"
$ws.Range("E10").Value = "'0.45%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07044"
$ws.Range("E11").Value = "'-0.58%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03203"
$ws.Range("E12").Value = "'9.54%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09283"
$ws.Range("E13").Value = "'-1.05%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001522"
$ws.Range("E14").Value = "'-0.05%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005993"
$ws.Range("E15").Value = "'-94.15%"
$ws.Range("D16").Value = "'0.006051"
$ws.Range("E16").Value = "'-1.74%"
$ws.Range("E17").Value = "'0.15%"
$ws.Range("D18").Value = "'2.175"
$ws.Range("E18").Value = "'-4.03%"
$ws.Range("E19").Value = "'-0.43%"
$ws.Range("D20").Value = "'0.03316"
$ws.Range("E20").Value = "'0.58%"
$ws.Range("E21").Value = "'-3.64%"
$ws.Range("D22").Value = "'3.499"
$ws.Range("E22").Value = "'1.01%"
$ws.Range("D23").Value = "'0.04093"
$ws.Range("E23").Value = "'-1.86%"
$ws.Range("D25").Value = "'0.001222"
$ws.Range("E25").Value = "'0.27%"
$ws.Range("D26").Value = "'0.004142"
$ws.Range("E26").Value = "'-17.70%"
$ws.Range("D27").Value = "'0.0001201"
$ws.Range("E27").Value = "'-0.79%"
$ws.Range("D28").Value = "'0.0001450"
$ws.Range("E28").Value = "'-25.21%"
$ws.Range("D40").Value = "'0.03755"
$ws.Range("E40").Value = "'0.28%"
$ws.Range("D41").Value = "'0.1063"
$ws.Range("E41").Value = "'-0.78%"
$ws.Range("D42").Value = "'0.003710"
$ws.Range("E42").Value = "'-36.47%"
$ws.Range("D43").Value = "'0.002449"
$ws.Range("E43").Value = "'22.46%"
$ws.Range("D44").Value = "'0.009365"
$ws.Range("E44").Value = "'-6.01%"
$ws.Range("D45").Value = "'0.00005263"
$ws.Range("E45").Value = "'1.15%"
$ws.Range("E46").Value = "'0.05%"
$ws.Range("D47").Value = "'0.07504"
$ws.Range("E47").Value = "'25.07%"
$ws.Range("D48").Value = "'0.002443"
$ws.Range("E48").Value = "'-4.71%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.05%"

Write-Host "Applied crypto price/volume updates"
